$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    11 = "1352d9b99bf06626ff80952eda02d7d2"
    34 = "c61e0c5fa0c3d3aeb7f195c62229f494"
    44 = "a2cfcbfef9b7b4aed5ed06cdf76e820f"
    74 = "9555bf74da8a390313ded720eb47dce7"
    89 = "160ee88f449d69ffbf488ebe9d2dcc44"
    99 = "ec5bd2a050b8a245967e920be6cdaaa2"
    110 = "4050bd447a74401c61ea746f9711d4fc"
    121 = "27c1bb70cb640d5ca20a759347c927c8"
    154 = "e9828e955ed4896624069e2230da5da2"
    160 = "f3de5288eeaf606f566c40f38f1f948a"
    161 = "9bb4c7968671c6ffbee5b3db18131f17"
    162 = "28b7081ddd8b2bf574091a34d8703cef"
    168 = "36c8cd53ba8a46717318adc0a51706b1"
    180 = "4452182d4a3e39871668d09fdb6c1e5b"
    191 = "c73e5ad0a567948972aa3db3a087d497"
    213 = "e11742ebab986b101aaf472dd8371e81"
    278 = "4f4e6e1d7f91885a3a4f184b8ac396e3"
    293 = "21201fdc44ce87e98d9209da669acf6b"
    335 = "ecbe729ac86df7acbe5e7934836f2f14"
    345 = "183913fecc02620ae6913e0667b17656"
    461 = "b11b80ec3b93464d6b97a5f9c1948435"
    480 = "f23b3dca7b162c63f81a3379142179f4"
    506 = "51d94fbb108c060af0774f3dfc25fd2e"
    514 = "1ff4dd27e25e4cecffa8c888a063c5c2"
    524 = "586802b4d9ba45de50d961c63708f3c0"
    534 = "76da3783aa2a61aa6867b6ba825b3179"
    547 = "12134a6651c6de21c72dc6c1e1dae89a"
    553 = "58d85ba2051dd71507a5e4255d2e5b94"
    572 = "f1eff8d1240251c266d684e4cbc1fca7"
    584 = "a576e1b2662d1a21d6c1d37626fd4452"
    666 = "6a504f8d367e29df8fe91b6e061f2350"
    729 = "27ed38bf1fbffac7273df8279ccba7ca"
    768 = "8a866f38cea4d509d812189b47eef642"
    811 = "dbd952bba9bedbb15ced3d14a76bc9b0"
    815 = "bd5b9380588c9dc7c9ba8123dc3cab76"
    816 = "1951623ae9020a139ec3467817acc2ab"
    825 = "76fb08e3968f1341beee8c4d704ab1a6"
    827 = "fe391b223dd9b3e7fc6a5f6ebd9890a3"
    874 = "d878f735a89572d2273c1e98708e28dd"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item([int]$row, 2).Value = $updates[$row]
}
